$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing data one column to the right
$ws.Columns("A").Insert()

# New section label: "N/A (split row and col)"
$ws.Range("B11").Value = "N/A (split row and col)"

# New "Threads" row header + thread counts (4, 3, 2)
$ws.Range("A12").Value = "Threads"
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 2

# New timing data rows
$ws.Range("B13").Value = 670978
$ws.Range("C13").Value = 809590
$ws.Range("D13").Value = 1030280

$ws.Range("B14").Value = 1937770
$ws.Range("C14").Value = 1937770
$ws.Range("D14").Value = 1937770

$ws.Range("B15").Formula = "=B14/B13"
$ws.Range("C15").Formula = "=C14/C13"
$ws.Range("D15").Formula = "=D14/D13"

# Apply styles matching existing style classes already used elsewhere in the sheet
# (references below use the POST-shift cell locations)
# style 1 (Menlo font, General number format) -> as used at C2:D2
$ws.Range("C2:D2").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)  # xlPasteFormats

# style 2 (numFmtId=11 'General number') -> as used at D3
$ws.Range("D3").Copy()
$ws.Range("D13").PasteSpecial(-4122)

# style 3 (numFmtId=11, Menlo font) -> as used at C7:D7
$ws.Range("C7:D7").Copy()
$ws.Range("B14:C14").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D14").PasteSpecial(-4122)

# style 2 for the new SpeedUp formula row
$ws.Range("D3").Copy()
$ws.Range("B15:D15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New column width for column G (closest achievable quantization to source 23.1640625)
$ws.Columns("G").ColumnWidth = 22.25

# Update the active selection to match the edited workbook state
$ws.Range("E12").Select() | Out-Null
